$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("A3").Value = 112379126
$ws.Range("B3").Value = 90826
$ws.Range("Q3").Value = 532036
$ws.Range("R3").Value = 6622648

# Row 4 updates
$ws.Range("A4").Value = 112379125
$ws.Range("B4").Value = 90826
$ws.Range("Q4").Value = 531963
$ws.Range("R4").Value = 6622561
